$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 3 (NamedFeatureSelector entry) - shifts rows 4,5,6 up to 3,4,5
$ws.Rows(3).Delete()

# Add new row 6 with style copied from row 5 (matches A column bordered/bold style)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value2 = 0
$ws.Range("B2").Value2 = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced',
                                    l1_ratio=0.7830765826716157, max_iter=1000,
                                    penalty='elasticnet', random_state=42,
                                    solver='saga'))])
'@
$ws.Range("C2").Value2 = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__l1_ratio': 0.7830765826716157, 'model__penalty': 'elasticnet', 'model__solver': 'saga'}
'@
$ws.Range("D2").Value2 = 0.570452602093469
$ws.Range("E2").Value2 = @'
Tree-Parzen Estimator
'@
$ws.Range("F2").Value2 = 42
$ws.Range("G2").Value2 = 0.7333838957458397
$ws.Range("H2").Value2 = 0.5440917107583775
$ws.Range("I2").Value2 = @'
[1 0 1 0 0 1 1 1 1 1 1 1 1 0 1 0 0 0 1 0 1 1 0 0]
'@
$ws.Range("J2").Value2 = @'
[0 1 1 1 0 0 1 0 0 1 1 1 1 1 1 0 1 0 0 0 1 0 0 1]
'@

# --- Row 3 ---
$ws.Range("A3").Value2 = 0
$ws.Range("B3").Value2 = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced',
                                    l1_ratio=0.8546943368620347, max_iter=1000,
                                    penalty='elasticnet', random_state=42,
                                    solver='saga'))])
'@
$ws.Range("C3").Value2 = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__l1_ratio': 0.8546943368620347, 'model__penalty': 'elasticnet', 'model__solver': 'saga'}
'@
$ws.Range("D3").Value2 = 0.5492912212000919
$ws.Range("E3").Value2 = @'
Tree-Parzen Estimator
'@
$ws.Range("F3").Value2 = 69
$ws.Range("G3").Value2 = 0.7439326998729044
$ws.Range("H3").Value2 = 0.5440917107583775
$ws.Range("I3").Value2 = @'
[0 1 1 0 1 0 0 0 1 1 1 0 1 0 1 0 1 0 1 1 0 1 1 1]
'@
$ws.Range("J3").Value2 = @'
[1 1 1 1 0 1 0 0 1 0 0 1 1 1 0 0 1 0 1 0 0 0 1 1]
'@

# --- Row 4 ---
$ws.Range("A4").Value2 = 0
$ws.Range("B4").Value2 = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced',
                                    l1_ratio=0.2732617604953749, max_iter=1000,
                                    penalty='elasticnet', random_state=42,
                                    solver='saga'))])
'@
$ws.Range("C4").Value2 = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__l1_ratio': 0.2732617604953749, 'model__penalty': 'elasticnet', 'model__solver': 'saga'}
'@
$ws.Range("D4").Value2 = 0.4826469469642844
$ws.Range("E4").Value2 = @'
Tree-Parzen Estimator
'@
$ws.Range("F4").Value2 = 23
$ws.Range("G4").Value2 = 0.6903718875045041
$ws.Range("H4").Value2 = 0.6666666666666666
$ws.Range("I4").Value2 = @'
[0 0 1 0 0 1 0 1 1 1 1 1 1 1 1 0 0 0 1 0 1 1 1 0]
'@
$ws.Range("J4").Value2 = @'
[0 0 1 0 1 1 1 1 1 1 1 0 1 1 1 1 0 1 0 0 1 0 0 0]
'@

# --- Row 5 ---
$ws.Range("A5").Value2 = 0
$ws.Range("B5").Value2 = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced',
                                    l1_ratio=0.4816613170749568, max_iter=1000,
                                    penalty='elasticnet', random_state=42,
                                    solver='saga'))])
'@
$ws.Range("C5").Value2 = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__l1_ratio': 0.4816613170749568, 'model__penalty': 'elasticnet', 'model__solver': 'saga'}
'@
$ws.Range("D5").Value2 = 0.5457854354556468
$ws.Range("E5").Value2 = @'
Tree-Parzen Estimator
'@
$ws.Range("F5").Value2 = 99
$ws.Range("G5").Value2 = 0.6793974437005478
$ws.Range("H5").Value2 = 0.4611594202898551
$ws.Range("I5").Value2 = @'
[0 0 1 1 0 1 1 1 1 0 1 1 0 1 1 0 0 1 0 1 1 0 0 1]
'@
$ws.Range("J5").Value2 = @'
[0 0 1 0 1 0 1 1 0 1 0 1 1 0 1 0 0 0 1 0 1 1 0 0]
'@

# --- Row 6 ---
$ws.Range("A6").Value2 = 0
$ws.Range("B6").Value2 = @'
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', None),
                ('model',
                 LogisticRegression(class_weight='balanced',
                                    l1_ratio=0.023731188293710237,
                                    max_iter=1000, penalty='elasticnet',
                                    random_state=42, solver='saga'))])
'@
$ws.Range("C6").Value2 = @'
{'selector': None, 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__l1_ratio': 0.023731188293710237, 'model__penalty': 'elasticnet', 'model__solver': 'saga'}
'@
$ws.Range("D6").Value2 = 0.548528176628863
$ws.Range("E6").Value2 = @'
Tree-Parzen Estimator
'@
$ws.Range("F6").Value2 = 89
$ws.Range("G6").Value2 = 0.6887664964018915
$ws.Range("H6").Value2 = 0.6269841269841269
$ws.Range("I6").Value2 = @'
[1 0 1 0 1 1 0 0 1 1 0 1 0 1 1 1 1 1 0 1 0 0 1 0]
'@
$ws.Range("J6").Value2 = @'
[0 0 1 1 1 0 0 0 1 1 1 1 0 1 1 0 0 0 0 1 0 1 1 1]
'@
